$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column D ("MSPE").
# This pushes the existing MSPE / Nugget / Non-expl columns one slot to
# the right (D->E, E->F, F->G), carrying their header style along.
$ws.Range("D1").EntireColumn.Insert()

# --- Header row -------------------------------------------------------
$ws.Range("D1").Value = "VAR(OK)"
$ws.Range("E1").Value = "MSPE"
$ws.Range("F1").Value = "S_nugget"
$ws.Range("G1").Value = "VAR(TOTAL)"
$ws.Range("H1").Value = "VAR(DATA)"

# Give the two brand-new header cells (G1, H1) the same bold / bordered
# / centered look as the rest of the header row (copy the formatting
# from an existing header cell rather than the brand-new empty style).
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# --- New VAR(OK) column (D2:D10) --------------------------------------
$ws.Range("D2").Value = 0.7558876508733949
$ws.Range("D3").Value = 1.403294474986249
$ws.Range("D4").Value = 1.699271575752888
$ws.Range("D5").Value = 1.791544649509852
$ws.Range("D6").Value = 1.868525521621687
$ws.Range("D7").Value = 1.910927781629038
$ws.Range("D8").Value = 1.9719854500803
$ws.Range("D9").Value = 1.982609598970315
$ws.Range("D10").Value = 1.987836712552757

# --- Column G used to hold the old "Non-expl var of model" values;
#     that metric no longer exists, so clear it back out (VAR(TOTAL)
#     stays blank for every row).
$ws.Range("G2:G10").ClearContents()

# --- New VAR(DATA) column (H2:H10), constant across all rows ----------
$ws.Range("H2:H10").Value = 8.40789956002731

# Dimension now spans A1:H10, matching the widened table.
